# Auto-applied numeric updates to Kujata_Profits market-data columns (H-N)
# across sheets ALC, ARM, BSM, CRP, GSM, LTW, WVR - scheduled data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 163.5
$ws.Range("I33").Value = 119.69231
$ws.Range("J33").Value = 353.33334
$ws.Range("K33").Value = 119.69231
$ws.Range("L33").Value = 353.33334
$ws.Range("M33").Value = 109.30769
$ws.Range("N33").Value = -811.33334

$ws.Range("H40").Value = 1931.4736
$ws.Range("J40").Value = 1546
$ws.Range("L40").Value = 1546
$ws.Range("N40").Value = -1896

$ws.Range("H76").Value = 5965.8887
$ws.Range("I76").Value = 5526.5
$ws.Range("J76").Value = 6091.4287
$ws.Range("K76").Value = 5526.5
$ws.Range("L76").Value = 6091.4287
$ws.Range("M76").Value = -5211.5
$ws.Range("N76").Value = -6721.4287

$ws.Range("H79").Value = 5965.8887
$ws.Range("I79").Value = 5526.5
$ws.Range("J79").Value = 6091.4287
$ws.Range("K79").Value = 5526.5
$ws.Range("L79").Value = 6091.4287
$ws.Range("M79").Value = -4434.5
$ws.Range("N79").Value = -8275.4287

$ws.Range("H80").Value = 388.22726
$ws.Range("J80").Value = 587.0833
$ws.Range("L80").Value = 1761.2499
$ws.Range("N80").Value = -3757.2499

$ws.Range("H83").Value = 388.22726
$ws.Range("J83").Value = 587.0833
$ws.Range("L83").Value = 5283.7497
$ws.Range("N83").Value = -15267.7497

$ws.Range("H86").Value = 4939.8945
$ws.Range("I86").Value = 3962.3333
$ws.Range("K86").Value = 3962.3333
$ws.Range("M86").Value = -2839.3333

$ws.Range("H88").Value = 1124661.5
$ws.Range("I88").Value = 602.5
$ws.Range("J88").Value = 1374452.4
$ws.Range("K88").Value = 602.5
$ws.Range("L88").Value = 1374452.4
$ws.Range("M88").Value = -196.5
$ws.Range("N88").Value = -1375264.4

$ws.Range("H89").Value = 4939.8945
$ws.Range("I89").Value = 3962.3333
$ws.Range("K89").Value = 19811.6665
$ws.Range("M89").Value = -14195.6665

$ws.Range("H91").Value = 1124661.5
$ws.Range("I91").Value = 602.5
$ws.Range("J91").Value = 1374452.4
$ws.Range("K91").Value = 602.5
$ws.Range("L91").Value = 1374452.4
$ws.Range("M91").Value = 801.5
$ws.Range("N91").Value = -1377260.4

$ws.Range("H116").Value = 3110.8286
$ws.Range("I116").Value = 2812.4092
$ws.Range("J116").Value = 3615.8462
$ws.Range("K116").Value = 2812.4092
$ws.Range("L116").Value = 3615.8462
$ws.Range("M116").Value = 629.5907999999999
$ws.Range("N116").Value = -10499.8462

$ws.Range("H138").Value = 464082.56
$ws.Range("I138").Value = 1213.6666
$ws.Range("J138").Value = 622780.5
$ws.Range("K138").Value = 3640.9998
$ws.Range("L138").Value = 1868341.5
$ws.Range("M138").Value = 1499.0002
$ws.Range("N138").Value = -1878621.5

$ws.Range("H141").Value = 2990
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 2990
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 8970
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -19330


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3777.239
$ws.Range("I32").Value = 3389.279
$ws.Range("K32").Value = 3389.279
$ws.Range("M32").Value = -3102.279

$ws.Range("H45").Value = 1067.5
$ws.Range("I45").Value = 1086.5217
$ws.Range("J45").Value = 980
$ws.Range("K45").Value = 1086.5217
$ws.Range("L45").Value = 980
$ws.Range("M45").Value = -709.5217
$ws.Range("N45").Value = -1734

$ws.Range("H61").Value = 1106.6
$ws.Range("I61").Value = 866.5714
$ws.Range("J61").Value = 1666.6666
$ws.Range("K61").Value = 866.5714
$ws.Range("L61").Value = 1666.6666
$ws.Range("M61").Value = -654.5714
$ws.Range("N61").Value = -2090.6666

$ws.Range("H74").Value = 1033.9697
$ws.Range("I74").Value = 565.04346
$ws.Range("J74").Value = 2112.5
$ws.Range("K74").Value = 565.04346
$ws.Range("L74").Value = 2112.5
$ws.Range("M74").Value = 308.95654
$ws.Range("N74").Value = -3860.5

$ws.Range("H77").Value = 1033.9697
$ws.Range("I77").Value = 565.04346
$ws.Range("J77").Value = 2112.5
$ws.Range("K77").Value = 2825.2173
$ws.Range("L77").Value = 10562.5
$ws.Range("M77").Value = 1542.7827
$ws.Range("N77").Value = -19298.5

$ws.Range("H132").Value = 1632.3334
$ws.Range("I132").Value = 1342.0513
$ws.Range("J132").Value = 2890.2222
$ws.Range("K132").Value = 4026.1539
$ws.Range("L132").Value = 8670.6666
$ws.Range("M132").Value = -1496.1539
$ws.Range("N132").Value = -13730.6666

$ws.Range("H136").Value = 1106.6
$ws.Range("I136").Value = 866.5714
$ws.Range("J136").Value = 1666.6666
$ws.Range("K136").Value = 2599.7142
$ws.Range("L136").Value = 4999.9998
$ws.Range("M136").Value = -49.71420000000035
$ws.Range("N136").Value = -10099.9998


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3070
$ws.Range("I20").Value = 2975
$ws.Range("J20").Value = 3133.3333
$ws.Range("K20").Value = 2975
$ws.Range("L20").Value = 3133.3333
$ws.Range("M20").Value = -2728
$ws.Range("N20").Value = -3627.3333


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 522.2727
$ws.Range("I22").Value = 343.125
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 343.125
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = 6.875
$ws.Range("N22").Value = -1700

$ws.Range("H29").Value = 10000
$ws.Range("J29").Value = 10000
$ws.Range("L29").Value = 10000
$ws.Range("N29").Value = -10586

$ws.Range("H31").Value = 1736.8334
$ws.Range("J31").Value = 2384.75
$ws.Range("L31").Value = 2384.75
$ws.Range("N31").Value = -2974.75

$ws.Range("H34").Value = 1736.8334
$ws.Range("J34").Value = 2384.75
$ws.Range("L34").Value = 2384.75
$ws.Range("N34").Value = -2788.75

$ws.Range("H107").Value = 637.5625
$ws.Range("I107").Value = 500.27274
$ws.Range("J107").Value = 939.6
$ws.Range("K107").Value = 500.27274
$ws.Range("L107").Value = 939.6
$ws.Range("M107").Value = 1419.72726
$ws.Range("N107").Value = -4779.6

$ws.Range("H122").Value = 807.8421
$ws.Range("I122").Value = 811.61536
$ws.Range("J122").Value = 799.6667
$ws.Range("K122").Value = 2434.84608
$ws.Range("L122").Value = 2399.0001
$ws.Range("M122").Value = 15.15391999999974
$ws.Range("N122").Value = -7299.0001

$ws.Range("H134").Value = 1619.4324
$ws.Range("I134").Value = 1674.6296
$ws.Range("J134").Value = 1470.4
$ws.Range("K134").Value = 5023.8888
$ws.Range("L134").Value = 4411.200000000001
$ws.Range("M134").Value = -2488.8888
$ws.Range("N134").Value = -9481.200000000001


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 56253224
$ws.Range("I70").Value = 50003420
$ws.Range("K70").Value = 50003420
$ws.Range("M70").Value = -50003150

$ws.Range("H73").Value = 56253224
$ws.Range("I73").Value = 50003420
$ws.Range("K73").Value = 50003420
$ws.Range("M73").Value = -50002484

$ws.Range("H132").Value = 1959.375
$ws.Range("I132").Value = 1796.1666
$ws.Range("J132").Value = 2449
$ws.Range("K132").Value = 5388.4998
$ws.Range("L132").Value = 7347
$ws.Range("M132").Value = -2858.4998
$ws.Range("N132").Value = -12407


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3767.8484
$ws.Range("I136").Value = 4512.885
$ws.Range("J136").Value = 1000.5714
$ws.Range("K136").Value = 13538.655
$ws.Range("L136").Value = 3001.7142
$ws.Range("M136").Value = -10988.655
$ws.Range("N136").Value = -8101.7142


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 166674670
$ws.Range("I62").Value = 250010000
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 250010000
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -250009376
$ws.Range("N62").Value = -5248

$ws.Range("H65").Value = 166674670
$ws.Range("I65").Value = 250010000
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 1250050000
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -1250046880
$ws.Range("N65").Value = -26240
